# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (holding the per-fund holdings detail)
#    right after "2021-Q4" and before "总计".
# 2) Insert a new summary row at the top of "总计" for the 2022-Q1 quarter,
#    pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: add the "2022-Q1" worksheet with the fund holdings detail table
# ---------------------------------------------------------------------------

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$newSheet.Name = "2022-Q1"
$ws = $newSheet

# existing "总计" sheet - used purely as a style donor so the new sheet reuses
# the workbook's existing header / index-column style instead of minting a new one
$totalSheet = $wb.Worksheets.Item("总计")

# sheetPr outline flags (summaryBelow="1" summaryRight="1")
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = -4152

# page margins -> 0.75in/0.75in/1in/1in/0.5in/0.5in (COM margins are in points, 72pt = 1in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Copy the header/index style (bold + border, centered) from the "总计" sheet
# onto the header row and the index column so we reuse the existing style slot.
$totalSheet.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Columns B:G hold numeric-looking codes/figures that are stored as TEXT
# (e.g. fund code "009387" must keep its leading zero). Format as text first
# so COM doesn't silently coerce them to numbers.
$ws.Range("B2:G6").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "009387"
$ws.Range("C2").Value = "嘉实稳福混合A"
$ws.Range("D2").Value = "0.08"
$ws.Range("E2").Value = "34.71"
$ws.Range("F2").Value = "4.26"
$ws.Range("G2").Value = "0.0034"
$ws.Range("H2").Value = 3

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "009649"
$ws.Range("C3").Value = "嘉实精选平衡混合A"
$ws.Range("D3").Value = "0.06"
$ws.Range("E3").Value = "67.70"
$ws.Range("F3").Value = "3.74"
$ws.Range("G3").Value = "0.0022"
$ws.Range("H3").Value = 8

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "001978"
$ws.Range("C4").Value = "泰信互联网+主题灵活配置混合"
$ws.Range("D4").Value = "0.06"
$ws.Range("E4").Value = "92.34"
$ws.Range("F4").Value = "2.60"
$ws.Range("G4").Value = "0.0016"
$ws.Range("H4").Value = 10

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "009650"
$ws.Range("C5").Value = "嘉实精选平衡混合C"
$ws.Range("D5").Value = "0.01"
$ws.Range("E5").Value = "67.70"
$ws.Range("F5").Value = "3.74"
$ws.Range("G5").Value = "0.0004"
$ws.Range("H5").Value = 8

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "009388"
$ws.Range("C6").Value = "嘉实稳福混合C"
$ws.Range("D6").Value = "0.01"
$ws.Range("E6").Value = "34.71"
$ws.Range("F6").Value = "4.26"
$ws.Range("G6").Value = "0.0004"
$ws.Range("H6").Value = 3

# Drop the temporary "@" text-format style again (value stays text, style goes
# back to the default/"Normal" slot, matching the un-styled detail cells used
# elsewhere in this workbook).
$ws.Range("B2:G6").Style = "Normal"

# ---------------------------------------------------------------------------
# Part 2: prepend a 2022-Q1 summary row onto the "总计" sheet
# ---------------------------------------------------------------------------

$tw = $totalSheet

# Shift the existing 3 data rows down by one (row4->row5, row3->row4, row2->row3),
# carrying the index-column style along with them.
$tw.Range("A4").Copy()
$tw.Range("A5").PasteSpecial(-4122)
$tw.Range("A5").Value = $tw.Range("A4").Value()
$tw.Range("B5").Value = $tw.Range("B4").Value()
$tw.Range("C5").Value = $tw.Range("C4").Value()
$tw.Range("D5").Value = $tw.Range("D4").Value()

$tw.Range("A4").Value = $tw.Range("A3").Value()
$tw.Range("B4").Value = $tw.Range("B3").Value()
$tw.Range("C4").Value = $tw.Range("C3").Value()
$tw.Range("D4").Value = $tw.Range("D3").Value()

$tw.Range("A3").Value = $tw.Range("A2").Value()
$tw.Range("B3").Value = $tw.Range("B2").Value()
$tw.Range("C3").Value = $tw.Range("C2").Value()
$tw.Range("D3").Value = $tw.Range("D2").Value()

# New top row: 2022-Q1 summary
$tw.Range("A2").Value = 0
$tw.Range("B2").Value = "2022-Q1"
$tw.Range("C2").Value = 5
$tw.Range("D2").Value = 0.01

# Restore the original active sheet/selection (the newly added sheet would
# otherwise become the active tab).
$wb.Worksheets.Item("2021-Q1").Activate()
